# Re-apply the refreshed cryptocurrency quotes (price/volume columns) that the
# "Updated cryptos list" GitHub Action produced, plus a few rows whose coins
# were re-ordered/re-ranked (name+link swapped with the following row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while preserving its original "General"
# look (no explicit style). Many of the new Price values are plain decimal
# numbers (e.g. "239.78"); assigning them straight to .Value would make Excel
# auto-convert the text into a real number. The source data keeps these as
# plain text, so when needed we briefly force a Text format, assign the
# string, then restore the default "Normal" style so no stray number format
# is left attached to the cell.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '26.327.94'
$ws.Range('E2').Value = '  +3.67%  '
$ws.Range('D3').Value = '1.719.05'
$ws.Range('E3').Value = '  +3.32%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue $ws.Range('D5') '239.78'
$ws.Range('E5').Value = '  +1.64%  '
$ws.Range('E6').Value = '  +0.01%  '
Set-TextValue $ws.Range('D7') '0.4714'
$ws.Range('E7').Value = '  -1.30%  '
Set-TextValue $ws.Range('D8') '0.2626'
$ws.Range('E8').Value = '  +0.76%  '
Set-TextValue $ws.Range('D9') '0.06216'
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('D10').Value = '1.716.96'
$ws.Range('E10').Value = '  +3.18%  '
Set-TextValue $ws.Range('D11') '0.07071'
$ws.Range('E11').Value = '  -0.01%  '
Set-TextValue $ws.Range('D12') '15.23'
$ws.Range('E12').Value = '  +3.29%  '
Set-TextValue $ws.Range('D13') '0.5913'
$ws.Range('E13').Value = '  -0.08%  '
Set-TextValue $ws.Range('D14') '4.412'
$ws.Range('E14').Value = '  +0.66%  '
Set-TextValue $ws.Range('D15') '76.37'
$ws.Range('E15').Value = '  +2.68%  '
Set-TextValue $ws.Range('D17') '1.001'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '26.335.99'
$ws.Range('E18').Value = '  +3.68%  '
Set-TextValue $ws.Range('D19') '0.000006792'
$ws.Range('E19').Value = '  +0.37%  '
Set-TextValue $ws.Range('D20') '11.60'
$ws.Range('E20').Value = '  +1.66%  '
$ws.Range('D21').Value = '1.933.53'
$ws.Range('E21').Value = '  +3.22%  '
Set-TextValue $ws.Range('D22') '4.557'
$ws.Range('E22').Value = '  +2.68%  '
Set-TextValue $ws.Range('D23') '8.791'
$ws.Range('E23').Value = '  +1.61%  '
Set-TextValue $ws.Range('D24') '5.338'
$ws.Range('E24').Value = '  -0.07%  '
Set-TextValue $ws.Range('D25') '134.88'
$ws.Range('E25').Value = '  +0.94%  '
Set-TextValue $ws.Range('D26') '15.17'
$ws.Range('E26').Value = '  +0.80%  '
Set-TextValue $ws.Range('D27') '1.405'
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('E28').Value = '  +4.26%  '
Set-TextValue $ws.Range('D29') '106.75'
$ws.Range('E29').Value = '  +2.43%  '
Set-TextValue $ws.Range('D30') '4.027'
$ws.Range('E30').Value = '  +0.91%  '
Set-TextValue $ws.Range('D31') '3.693'
$ws.Range('E31').Value = '  +2.23%  '
Set-TextValue $ws.Range('D32') '0.07727'
$ws.Range('E32').Value = '  +1.11%  '
Set-TextValue $ws.Range('D33') '0.04442'
$ws.Range('E33').Value = '  +1.47%  '
Set-TextValue $ws.Range('D34') '2.613'
$ws.Range('E34').Value = '  +0.35%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D35') '0.9727'
$ws.Range('E35').Value = '  +3.05%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D36') '0.6199'
$ws.Range('E36').Value = '  +1.20%  '
$ws.Range('B37').Value = 'Quant'
$ws.Range('C37').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D37') '115.70'
$ws.Range('E37').Value = '  +17.76%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D38') '0.9291'
$ws.Range('E38').Value = '  +8.71%  '
Set-TextValue $ws.Range('D39') '2.412'
$ws.Range('E39').Value = '  -8.08%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range('D40') '1.001'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D41') '1.904'
$ws.Range('E41').Value = '  +4.03%  '
Set-TextValue $ws.Range('D42') '0.01471'
$ws.Range('E42').Value = '  -2.27%  '
Set-TextValue $ws.Range('D43') '5.297'
$ws.Range('E43').Value = '  +13.52%  '
$ws.Range('E44').Value = '  +1.35%  '
Set-TextValue $ws.Range('D45') '0.1156'
$ws.Range('E45').Value = '  +3.48%  '
Set-TextValue $ws.Range('D46') '6.258'
$ws.Range('E46').Value = '  +0.76%  '
$ws.Range('E47').Value = '  +0.66%  '
Set-TextValue $ws.Range('D48') '30.62'
$ws.Range('E48').Value = '  +3.63%  '
Set-TextValue $ws.Range('D49') '7.648'
$ws.Range('E49').Value = '  +4.87%  '
$ws.Range('E50').Value = '  +1.33%  '
Set-TextValue $ws.Range('D51') '1.219'
$ws.Range('E51').Value = '  +1.37%  '
